# Add a new column F ("Totales" as of 2020-06-11 / serial 43993) that
# mirrors column E's running-total table, with its own day-by-day values
# (mostly identical to column E, but a handful of the later rows carry a
# slightly higher count) and its own SUM() total in row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 header cell (date) -------------------------------------------
# Give F2 the same border as the rest of the table plus the yellow fill +
# "d-mmm" date format used by the other header date cells (B2:E2), but
# without forcing center alignment (matches the newly authored style).
$ws.Range("E3").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F2").NumberFormat = "d-mmm"
$ws.Range("F2").Interior.Color = 65535
$ws.Range("F2").Value = 43993

# --- Day-by-day values (rows 3-80) ---------------------------------------
# Bulk-copy the existing formatting from column E so the same style indexes
# get reused (border-only for rows 3-76, white-fill+border for rows 77-79).
$ws.Range("E3:E76").Copy()
$ws.Range("F3:F76").PasteSpecial(-4122)

$ws.Range("E77:E79").Copy()
$ws.Range("F77:F79").PasteSpecial(-4122)

$values = @(1,1,1,1,1,1,3,4,1,4,6,8,1,5,4,10,9,6,7,3,9,4,8,6,5,6,5,10,8,5,7,6,7,9,11,7,6,8,10,7,19,12,11,14,29,16,22,30,19,25,38,35,32,47,45,60,60,65,60,69,84,93,85,101,111,109,117,119,105,128,121,124,112,46,67,98,79)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 3, 6).Value = $values[$i]
}

# Row 80 keeps the plain (unformatted) default style - no format copied.
$ws.Range("F80").Value = 20

# --- Row 1 total ------------------------------------------------------
$ws.Range("F1").Formula = "=SUM(F3:F80)"

# --- Misc cosmetic refresh to match the re-saved file ---------------------
$excel.ActiveWindow.Zoom = 71
$ws.Range("R19").Select()
